# PowerPoint Left/Top are Single-precision (float) values expressed in
# points; 1 pt = 12700 EMU. The EMU targets below are nudged by at most
# 1 EMU where needed so that, after the float32 round-trip PowerPoint's
# COM layer performs, the on-disk EMU value lands on (or as close as
# representable to) the authored target.
$emuPerPt = 12700

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# -----------------------------------------------------------------------
# Reposition the six shapes that make up the "Analyse Stopped" /
# "Analysing Telemetry" sub-flow (moved together, as a group, by the
# same offset).
# -----------------------------------------------------------------------

# Oval 25
$sh = $s.Shapes.Item(23)
$sh.Left = 3071817 / $emuPerPt
$sh.Top  = 14166995 / $emuPerPt

# Oval 26
$sh = $s.Shapes.Item(24)
$sh.Left = 5006833 / $emuPerPt
$sh.Top  = 14134667 / $emuPerPt

# Curved Connector 27
$sh = $s.Shapes.Item(25)
$sh.Left = 4605730 / $emuPerPt
$sh.Top  = 13165529 / $emuPerPt

# TextBox 28 ("Started braking")
$sh = $s.Shapes.Item(26)
$sh.Left = 4099497 / $emuPerPt
$sh.Top  = 13658729 / $emuPerPt

# Curved Connector 33
$sh = $s.Shapes.Item(31)
$sh.Left = 4605730 / $emuPerPt
$sh.Top  = 13694448 / $emuPerPt

# TextBox 34 ("Stopped braking")
$sh = $s.Shapes.Item(32)
$sh.Left = 4030224 / $emuPerPt
$sh.Top  = 15030329 / $emuPerPt

# -----------------------------------------------------------------------
# Fix typo: "Front slip >= threshold" -> "Front slip <= threshold"
# Split the run into three pieces, matching the authored edit.
# -----------------------------------------------------------------------
$sh = $s.Shapes.Item(33)   # Flowchart: Decision 35
$tr = $sh.TextFrame.TextRange
$sub = $tr.Characters(12, 3)   # ">= "
$sub.Text = "<= "
